$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Doing Updates for Financials": two more quarters of data have come in for
# DIT (period endings 2018-12-31 and 2018-09-30), so insert two fresh columns
# right before the existing "most recent quarter" column (D) and shift all
# the historical quarters two columns to the right (D:K -> F:M).
# ---------------------------------------------------------------------------
$ws.Range("D1:E1").EntireColumn.Insert()

# The newly inserted D:E columns should carry the same per-row number
# formatting (dates vs plain numbers) that column D used to have before the
# insert -- which, after the shift above, now lives in column F. Copy formats
# across from F before any values are written into D/E.
$ws.Range("F1:F1048576").EntireColumn.Copy()
$ws.Range("D1:E1048576").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Keep the new columns the same width as their neighbours.
$ws.Range("D1").ColumnWidth = $ws.Range("F1").ColumnWidth
$ws.Range("E1").ColumnWidth = $ws.Range("F1").ColumnWidth

# Populate the two new quarterly columns: D = period ending 2018-12-31
# (serial 43465), E = period ending 2018-09-30 (serial 43373).
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(8, 4).Value = 344700
$ws.Cells.Item(8, 5).Value = 362500
$ws.Cells.Item(9, 4).Value = 324100
$ws.Cells.Item(9, 5).Value = 340000
$ws.Cells.Item(10, 4).Value = 20600
$ws.Cells.Item(10, 5).Value = 22500
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 4).Value = "NA"
$ws.Cells.Item(14, 5).Value = 1900
$ws.Cells.Item(15, 4).Value = 600
$ws.Cells.Item(15, 5).Value = 600
$ws.Cells.Item(17, 4).Value = 342700
$ws.Cells.Item(17, 5).Value = 360300
$ws.Cells.Item(18, 4).Value = 2000
$ws.Cells.Item(18, 5).Value = 2200
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(21, 4).Value = 2700
$ws.Cells.Item(21, 5).Value = 2900
$ws.Cells.Item(22, 4).Value = 300
$ws.Cells.Item(22, 5).Value = 400
$ws.Cells.Item(23, 4).Value = 1700
$ws.Cells.Item(23, 5).Value = 1800
$ws.Cells.Item(24, 4).Value = 1300
$ws.Cells.Item(24, 5).Value = 700
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = 400
$ws.Cells.Item(26, 5).Value = 1100
$ws.Cells.Item(27, 4).Value = 400
$ws.Cells.Item(27, 5).Value = 1100
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 4).Value = 800
$ws.Cells.Item(29, 5).Value = -100
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(33, 4).Value = 1200
$ws.Cells.Item(33, 5).Value = 1000
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 4).Value = 1200
$ws.Cells.Item(35, 5).Value = 1000
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(41, 4).Value = 800
$ws.Cells.Item(41, 5).Value = 500
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(43, 4).Value = 30800
$ws.Cells.Item(43, 5).Value = 31700
$ws.Cells.Item(44, 4).Value = 56600
$ws.Cells.Item(44, 5).Value = 78900
$ws.Cells.Item(45, 4).Value = 9300
$ws.Cells.Item(45, 5).Value = 4900
$ws.Cells.Item(46, 4).Value = 97400
$ws.Cells.Item(46, 5).Value = 116000
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(48, 4).Value = 16300
$ws.Cells.Item(48, 5).Value = 15800
$ws.Cells.Item(49, 4).Value = 7800
$ws.Cells.Item(49, 5).Value = 7900
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 4).Value = 300
$ws.Cells.Item(52, 5).Value = 300
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 4).Value = 121800
$ws.Cells.Item(54, 5).Value = 140000
$ws.Cells.Item(57, 4).Value = 16800
$ws.Cells.Item(57, 5).Value = 20800
$ws.Cells.Item(58, 4).Value = 1000
$ws.Cells.Item(58, 5).Value = 1100
$ws.Cells.Item(59, 4).Value = 9500
$ws.Cells.Item(59, 5).Value = 12500
$ws.Cells.Item(60, 4).Value = 27200
$ws.Cells.Item(60, 5).Value = 34400
$ws.Cells.Item(61, 4).Value = 26900
$ws.Cells.Item(61, 5).Value = 39100
$ws.Cells.Item(62, 4).Value = 2100
$ws.Cells.Item(62, 5).Value = 1800
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 4).Value = 56200
$ws.Cells.Item(66, 5).Value = 75400
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 4).Value = 64800
$ws.Cells.Item(72, 5).Value = 63800
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 4).Value = 65700
$ws.Cells.Item(76, 5).Value = 64600
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(81, 4).Value = 1200
$ws.Cells.Item(81, 5).Value = 1000
$ws.Cells.Item(83, 4).Value = 600
$ws.Cells.Item(83, 5).Value = 600
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 4).Value = 14600
$ws.Cells.Item(89, 5).Value = -13200
$ws.Cells.Item(91, 4).Value = -1000
$ws.Cells.Item(91, 5).Value = -1100
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 4).Value = -1000
$ws.Cells.Item(94, 5).Value = -1100
$ws.Cells.Item(96, 4).Value = -100
$ws.Cells.Item(96, 5).Value = -100
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 4).Value = -13400
$ws.Cells.Item(100, 5).Value = 14500
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(101, 5).Value = 0
$ws.Cells.Item(102, 4).Value = 300
$ws.Cells.Item(102, 5).Value = 200


# "Discontinued Operations" (row 14) reported actual figures for the two new
# quarters plus the next older one, but the three quarters after that
# switched from "0" to not-applicable ("NA") instead of the carried-over
# zeroes that the plain column shift would otherwise leave behind.
$ws.Cells.Item(14, 6).Value = "NA"
$ws.Cells.Item(14, 7).Value = "NA"
$ws.Cells.Item(14, 8).Value = "NA"
$ws.Cells.Item(14, 9).Value = "NA"
$ws.Cells.Item(14, 10).Value = "NA"
